# Update GDP Growth Rates workbook to the November 2020 STEO / EPS 3.1 refresh.

$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")
$data = $wb.Worksheets.Item("Data")

# --- "Data" sheet: swap in the November STEO vintage figures ---
$data.Range("A3").Value = "November STEO"
$data.Range("B3").Value = 19092
$data.Range("C3").Value = 18411
$data.Range("D3").Value = 19098

# --- "About" sheet: refresh narrative text referencing the STEO vintage/version ---
$about.Range("B6").Value = "January 2020 and November 2020"
$about.Range("A27").Value = "As of EPS 3.1, this variable is set up to model the impacts of the 2020"
$about.Range("A28").Value = "SARS-CoV-2 pandemic.  It uses the latest data available as of November 10,"

# Restore the saved selection state to match the authored workbook, without
# changing which sheet is active (the "About" sheet stays the active tab).
$originalActive = $wb.ActiveSheet
$data.Activate()
$data.Range("B12").Select()
$originalActive.Activate()
